$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new user record (row 33) mirroring the existing data pattern
$row = 33
$ws.Cells.Item($row, 1).Value = 110032
$ws.Cells.Item($row, 2).Value = 9317596770
$ws.Cells.Item($row, 3).Value = "Ewan Marsh"
$ws.Cells.Item($row, 4).Value = "ewan.marsh@xyz.com"
$ws.Cells.Item($row, 5).Value = 818876433
$ws.Cells.Item($row, 6).Value = "ACT"
$ws.Cells.Item($row, 7).Value = "eng"
$ws.Cells.Item($row, 8).Value = "PWD"
$ws.Cells.Item($row, 9).Value = $true
$ws.Cells.Item($row, 10).Value = "superadmin"
$ws.Cells.Item($row, 11).Value = "now()"

# Match the left-aligned style already used for the "email" and
# "is_active" columns in the existing rows
$ws.Cells.Item($row, 4).HorizontalAlignment = -4131
$ws.Cells.Item($row, 9).HorizontalAlignment = -4131

# Scroll back to the top and select the full column L (beyond the data range)
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("L1:XFD1048576").Select()

# Restore print setup vertical DPI
$ws.PageSetup.VerticalDpi = 300

$wb.Save()
